# TournRPG-138: update enemy parameter table (exp requirements / lv1 hp)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("player")

# New "exp" (column B) values for rows 3..102 (lv 1..100)
$expValues = @(
    10, 25, 100, 200, 300, 400, 500, 630, 760, 900,
    1050, 1200, 1350, 1500, 1700, 2000, 2500, 3000, 3500, 4000,
    4500, 5000, 5500, 6000, 6500, 7000, 7500, 8000, 8500, 9000,
    9500, 10000, 10500, 11000, 11500, 12000, 12500, 13000, 13500, 14000,
    14500, 15000, 15500, 16000, 16500, 17000, 17500, 18000, 18500, 19000,
    19500, 20000, 20500, 21000, 21500, 22000, 22500, 23000, 23500, 24000,
    24500, 25000, 25500, 26000, 26500, 27000, 27500, 28000, 28500, 29000,
    29500, 30000, 30500, 31000, 31500, 32000, 32500, 33000, 33500, 34000,
    34500, 35000, 35500, 36000, 36500, 37000, 37500, 38000, 38500, 39000,
    39500, 40000, 40500, 41000, 41500, 42000, 42500, 43000, 43500, 44000
)

$startRow = 3
for ($i = 0; $i -lt $expValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $expValues[$i]
}

# hp (column C) for lv1 (row 3) changes from 100 to 50
$ws.Cells.Item(3, 3).Value = 50
